$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated state-transition probabilities after adding more simulated games
# (recomputed frequencies from the larger sample of games), per the
# "added more games" commit.
$ws.Range("B2").Value = 0.1956521739130435
$ws.Range("C2").Value = 0.5434782608695652
$ws.Range("J2").Value = 0.01304347826086956
$ws.Range("P2").Value = 0.108695652173913
$ws.Range("S2").Value = 0.1391304347826087
$ws.Range("B3").Value = 0.007936507936507936
$ws.Range("C3").Value = 0.007936507936507936
$ws.Range("J3").Value = 0.01587301587301587
$ws.Range("P3").Value = 0.753968253968254
$ws.Range("S3").Value = 0.2142857142857143
$ws.Range("J4").Value = 0.02439024390243903
$ws.Range("P4").Value = 0.5121951219512195
$ws.Range("S4").Value = 0.4634146341463415
$ws.Range("B6").Value = 0.04020100502512563
$ws.Range("D6").Value = 0.01507537688442211
$ws.Range("E6").Value = 0.005025125628140704
$ws.Range("F6").Value = 0.04020100502512563
$ws.Range("J6").Value = 0.2311557788944724
$ws.Range("O6").Value = 0.02512562814070352
$ws.Range("Q6").Value = 0.1708542713567839
$ws.Range("R6").Value = 0.06532663316582915
$ws.Range("S6").Value = 0.407035175879397
$ws.Range("B7").Value = 0.1217948717948718
$ws.Range("D7").Value = 0.00641025641025641
$ws.Range("F7").Value = 0.0641025641025641
$ws.Range("J7").Value = 0.09615384615384616
$ws.Range("O7").Value = 0.02564102564102564
$ws.Range("R7").Value = 0.07051282051282051
$ws.Range("S7").Value = 0.4615384615384616
$ws.Range("B8").Value = 0.06167400881057269
$ws.Range("D8").Value = 0.006607929515418502
$ws.Range("E8").Value = 0.002202643171806168
$ws.Range("F8").Value = 0.05506607929515418
$ws.Range("J8").Value = 0.1167400881057269
$ws.Range("O8").Value = 0.013215859030837
$ws.Range("Q8").Value = 0.1784140969162996
$ws.Range("R8").Value = 0.07268722466960352
$ws.Range("S8").Value = 0.4933920704845815
$ws.Range("B9").Value = 0.06278026905829596
$ws.Range("D9").Value = 0.01345291479820628
$ws.Range("E9").Value = 0.004484304932735426
$ws.Range("F9").Value = 0.07623318385650224
$ws.Range("J9").Value = 0.1031390134529148
$ws.Range("O9").Value = 0.01345291479820628
$ws.Range("Q9").Value = 0.1838565022421525
$ws.Range("R9").Value = 0.07623318385650224
$ws.Range("S9").Value = 0.4663677130044843
$ws.Range("B10").Value = 0.09105960264900662
$ws.Range("D10").Value = 0.02649006622516556
$ws.Range("E10").Value = 0.0008278145695364238
$ws.Range("F10").Value = 0.06622516556291391
$ws.Range("J10").Value = 0.1142384105960265
$ws.Range("O10").Value = 0.009933774834437087
$ws.Range("Q10").Value = 0.1928807947019868
$ws.Range("R10").Value = 0.09437086092715231
$ws.Range("S10").Value = 0.4039735099337748
$ws.Range("G11").Value = 0.1611570247933884
$ws.Range("J11").Value = 0.1074380165289256
$ws.Range("K11").Value = 0.2272727272727273
$ws.Range("L11").Value = 0.487603305785124
$ws.Range("S11").Value = 0.01652892561983471
$ws.Range("G12").Value = 0.7166666666666667
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.008333333333333333
$ws.Range("L12").Value = 0.01666666666666667
$ws.Range("S12").Value = 0.05833333333333333
$ws.Range("G13").Value = 0.7906976744186046
$ws.Range("J13").Value = 0.2093023255813954
$ws.Range("F15").Value = 0.01036269430051814
$ws.Range("H15").Value = 0.1450777202072539
$ws.Range("I15").Value = 0.08290155440414508
$ws.Range("J15").Value = 0.4196891191709844
$ws.Range("K15").Value = 0.05699481865284974
$ws.Range("M15").Value = 0.02072538860103627
$ws.Range("N15").Value = 0.005181347150259068
$ws.Range("O15").Value = 0.05181347150259067
$ws.Range("S15").Value = 0.2072538860103627
$ws.Range("F16").Value = 0.01418439716312057
$ws.Range("H16").Value = 0.1843971631205674
$ws.Range("I16").Value = 0.09929078014184398
$ws.Range("J16").Value = 0.5106382978723404
$ws.Range("K16").Value = 0.07092198581560284
$ws.Range("M16").Value = 0.02127659574468085
$ws.Range("O16").Value = 0.0425531914893617
$ws.Range("S16").Value = 0.05673758865248227
$ws.Range("F17").Value = 0.009708737864077669
$ws.Range("H17").Value = 0.1966019417475728
$ws.Range("I17").Value = 0.08495145631067962
$ws.Range("J17").Value = 0.4320388349514563
$ws.Range("K17").Value = 0.06796116504854369
$ws.Range("M17").Value = 0.03883495145631068
$ws.Range("O17").Value = 0.0703883495145631
$ws.Range("S17").Value = 0.09951456310679611
$ws.Range("F18").Value = 0.0273224043715847
$ws.Range("H18").Value = 0.2021857923497268
$ws.Range("I18").Value = 0.1420765027322404
$ws.Range("J18").Value = 0.4043715846994536
$ws.Range("K18").Value = 0.06010928961748634
$ws.Range("M18").Value = 0.0273224043715847
$ws.Range("N18").Value = 0.00546448087431694
$ws.Range("O18").Value = 0.04918032786885246
$ws.Range("S18").Value = 0.08196721311475409
$ws.Range("F19").Value = 0.01471727343144849
$ws.Range("H19").Value = 0.2199845081332301
$ws.Range("I19").Value = 0.104570100697134
$ws.Range("J19").Value = 0.372579395817196
$ws.Range("K19").Value = 0.09527498063516654
$ws.Range("M19").Value = 0.01316808675445391
$ws.Range("N19").Value = 0.001549186676994578
$ws.Range("O19").Value = 0.06738962044926414
$ws.Range("S19").Value = 0.1107668474051123
